# 148. Sort List and 2095. Delete the middle Node of linked list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 - LeetCode 2095: Delete the Middle Node of a Linked List
$ws.Cells.Item(32, 1).Value = 2095
$ws.Cells.Item(32, 2).Value = "Delete the Middle Node of a Linked List"
$ws.Cells.Item(32, 3).Value = "Java"

# Row 33 - LeetCode 148: Sort List
$ws.Cells.Item(33, 1).Value = 148
$ws.Cells.Item(33, 2).Value = "Sort List"
$ws.Cells.Item(33, 3).Value = "Java"

# Update the view so the newly added rows are visible/selected,
# matching the author's saved cursor position.
$ws.Range("C33").Select()
